# Commit: "Changed delimeter from _ to ."
# Rename both sheets so the underscore separator becomes a dot separator,
# then restore each sheet's last active-cell selection.

$wb = $excel.ActiveWorkbook

$sheetRegions = $wb.Worksheets.Item("dcim_regions")
$sheetSites   = $wb.Worksheets.Item("dcim_sites")

$sheetRegions.Name = "dcim.regions"
$sheetSites.Name   = "dcim.sites"

# dcim.regions (first sheet) was left with G31 selected.
$sheetRegions.Activate()
$sheetRegions.Range("G31").Select()

# dcim.sites (second/active sheet) was left with F29 selected.
$sheetSites.Activate()
$sheetSites.Range("F29").Select()
